$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 12:50:24"
$wsZhCn.Range("H2").Value = "2016-03-12 12:50:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 12:50:27"
$wsDeDe.Range("H2").Value = "2016-03-12 12:50:48"
